$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.739.63'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '2.454.17'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '570.94'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').Value = '146.03'
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -1.74%  '
$ws.Range('E9').Value = '  -1.06%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('D11').Value = '5.17'
$ws.Range('E11').Value = '  -2.22%  '
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('D13').Value = '28.70'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('E14').Value = '  -3.35%  '
$ws.Range('D15').Value = '2.900.25'
$ws.Range('D16').Value = '62.520.98'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('D17').Value = '2.416.51'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('D18').Value = '7.66'
$ws.Range('E18').Value = '  -5.66%  '
$ws.Range('D19').Value = '10.73'
$ws.Range('E19').Value = '  -3.11%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = '321.11'
$ws.Range('D22').Value = '2.21'
$ws.Range('E22').Value = '  -1.22%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '9.88'
$ws.Range('E24').Value = '  +3.04%  '
$ws.Range('D25').Value = '65.05'
$ws.Range('D26').Value = '644.31'
$ws.Range('E26').Value = '  -4.33%  '
$ws.Range('D27').Value = '2.574.49'
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('D28').Value = '0.0₃0955'
$ws.Range('E28').Value = '  -4.17%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('D31').Value = '7.81'
$ws.Range('E31').Value = '  -3.74%  '
$ws.Range('E32').Value = '  -3.46%  '
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  -4.11%  '
$ws.Range('E36').Value = '  -3.40%  '
$ws.Range('D37').Value = '151.71'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('D38').Value = '0.363'
$ws.Range('D39').Value = '18.49'
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('D40').Value = '5.30'
$ws.Range('E40').Value = '  -3.51%  '
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('E42').Value = '  -2.91%  '
$ws.Range('D43').Value = '0.0₆0308'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = '152.40'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').Value = '15.40'
$ws.Range('E46').Value = '  +1.50%  '
$ws.Range('E47').Value = '  -2.78%  '
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('D49').Value = '20.05'
$ws.Range('E49').Value = '  -3.96%  '
$ws.Range('D50').Value = '0.0502'
$ws.Range('E50').Value = '  -2.71%  '
$ws.Range('D51').Value = '0.0905'
$ws.Range('E51').Value = '  -1.67%  '
